# Updated symbol list on Mon Dec 26 02:54:41 UTC 2022 with GitHub Actions
#
# This script reproduces the data refresh captured in the target diff:
#  - several "Price" (column D) values are refreshed to newer quotes
#  - a new coin ("One") is inserted at row 10 (with a "Bestin24h" volume tag),
#    pushing WazirX..CoinExToken down by one row each (rows 10-18), with
#    each of those rows also getting refreshed Price values
#  - a couple of "Volume(1h)" (column E) labels lose/gain a "Bestin24h" suffix
#
# Column D in this sheet stores prices as plain text (not numbers), so we
# must force text storage via a leading apostrophe - otherwise Excel's
# automatic type inference would convert these into real floating point
# numbers (losing exact formatting such as trailing zeros).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $value) {
    # Force the cell to store an exact text value, even if it looks numeric.
    $ws.Range($cellRef).Value = "'" + $value
}

function Set-StringCell($cellRef, $value) {
    # Plain (non numeric-looking) text, safe to assign directly.
    $ws.Range($cellRef).Value = $value
}

# --- Price refreshes above the inserted row ---------------------------------
Set-TextCell "D2" "244.07"
Set-TextCell "D3" "23.21"
Set-TextCell "D4" "5.411"
Set-TextCell "D5" "0.05989"
Set-TextCell "D6" "3.464"
Set-TextCell "D7" "6.518"
Set-TextCell "D8" "0.8162"
Set-TextCell "D9" "0.9204"

# --- Row 10: "One" is newly inserted here (was previously at row 18) --------
Set-StringCell "B10" "One"
Set-StringCell "C10" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextCell   "D10" "0.01128"
Set-StringCell "E10" "9OneONEBestin24h"

# --- Row 11: WazirX (shifted down from row 10) -------------------------------
Set-StringCell "B11" "WazirX"
Set-StringCell "C11" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextCell   "D11" "0.1408"
Set-StringCell "E11" "10WazirXWRX"

# --- Row 12: MandalaExchangeToken (shifted down from row 11) ----------------
Set-StringCell "B12" "MandalaExchangeToken"
Set-StringCell "C12" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextCell   "D12" "0.07395"
Set-StringCell "E12" "11MandalaExchangeTokenMDX"

# --- Row 13: LiechtensteinCryptoassetsExchange (shifted down from row 12) ---
Set-StringCell "B13" "LiechtensteinCryptoassetsExchange"
Set-StringCell "C13" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextCell   "D13" "0.03247"
Set-StringCell "E13" "12LiechtensteinCryptoassetsExchangeLCX"

# --- Row 14: BitrueCoin (shifted down from row 13) --------------------------
Set-StringCell "B14" "BitrueCoin"
Set-StringCell "C14" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextCell   "D14" "0.03057"
Set-StringCell "E14" "13BitrueCoinBTR"

# --- Row 15: BitMartToken (shifted down from row 14) ------------------------
Set-StringCell "B15" "BitMartToken"
Set-StringCell "C15" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextCell   "D15" "0.09345"
Set-StringCell "E15" "14BitMartTokenBMX"

# --- Row 16: MCDex (shifted down from row 15) -------------------------------
Set-StringCell "B16" "MCDex"
Set-StringCell "C16" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextCell   "D16" "3.879"
Set-StringCell "E16" "15MCDexMCB"

# --- Row 17: BitForexToken (shifted down from row 16) -----------------------
Set-StringCell "B17" "BitForexToken"
Set-StringCell "C17" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextCell   "D17" "0.001559"
Set-StringCell "E17" "16BitForexTokenBF"

# --- Row 18: CoinExToken (shifted down from row 17) -------------------------
Set-StringCell "B18" "CoinExToken"
Set-StringCell "C18" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextCell   "D18" "0.04682"
Set-StringCell "E18" "17CoinExTokenCET"

# --- Remaining scattered price / label refreshes ----------------------------
Set-TextCell   "D19" "0.006097"

Set-TextCell   "D20" "0.005012"
Set-StringCell "E20" "19HotbitTokenHTB"

Set-TextCell "D21" "0.0009851"
Set-TextCell "D22" "0.00007800"
Set-TextCell "D24" "2.130"
Set-TextCell "D27" "0.0002900"
Set-TextCell "D41" "0.006217"
Set-TextCell "D42" "0.1074"
Set-TextCell "D43" "0.003000"
Set-TextCell "D44" "0.007134"
Set-TextCell "D45" "0.00005229"
Set-TextCell "D48" "0.9100"
Set-TextCell "D49" "0.002299"
Set-TextCell "D50" "0.00002100"
Set-TextCell "D51" "0.0002000"
